$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "30.197.07"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.910.64"
$ws.Range("E3").Value = "  -0.13%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "0.8234"
$ws.Range("E5").Value = "  +4.62%  "
Set-TextValue "D6" "243.74"
$ws.Range("E6").Value = "  +0.26%  "
Set-TextValue "D7" "1.003"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +2.86%  "
Set-TextValue "D9" "26.90"
$ws.Range("E9").Value = "  +2.40%  "
Set-TextValue "D10" "0.07059"
$ws.Range("E10").Value = "  +1.79%  "
Set-TextValue "D11" "0.08101"
$ws.Range("E11").Value = "  +1.39%  "
Set-TextValue "D12" "0.7665"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.911.18"
$ws.Range("E13").Value = "  -0.07%  "
Set-TextValue "D14" "5.285"
$ws.Range("E14").Value = "  +1.21%  "
Set-TextValue "D15" "92.90"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "30.199.99"
$ws.Range("E16").Value = "  +0.33%  "
Set-TextValue "D17" "14.21"
$ws.Range("E17").Value = "  +1.48%  "
Set-TextValue "D18" "5.902"
$ws.Range("E18").Value = "  -0.44%  "
Set-TextValue "D19" "245.47"
$ws.Range("E19").Value = "  -0.75%  "
Set-TextValue "D20" "0.000007794"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "2.163.42"
$ws.Range("E21").Value = "  -0.02%  "
Set-TextValue "D22" "1.004"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue "D24" "7.043"
$ws.Range("E24").Value = "  +1.95%  "
Set-TextValue "D25" "0.1676"
$ws.Range("E25").Value = "  +20.95%  "
Set-TextValue "D26" "9.317"
$ws.Range("E26").Value = "  +0.02%  "
Set-TextValue "D27" "166.78"
$ws.Range("E27").Value = "  -1.38%  "
Set-TextValue "D28" "19.02"
$ws.Range("E28").Value = "  +0.44%  "
Set-TextValue "D29" "2.108"
$ws.Range("E29").Value = "  +3.04%  "
Set-TextValue "D30" "1.373"
$ws.Range("E30").Value = "  -0.44%  "
Set-TextValue "D31" "1.526"
$ws.Range("E31").Value = "  +0.16%  "
Set-TextValue "D32" "0.05956"
$ws.Range("E32").Value = "  +3.91%  "
Set-TextValue "D33" "4.306"
$ws.Range("E33").Value = "  -0.87%  "
Set-TextValue "D34" "4.092"
$ws.Range("E34").Value = "  -0.59%  "
Set-TextValue "D35" "1.274"
$ws.Range("E35").Value = "  +0.94%  "
Set-TextValue "D36" "0.7346"
$ws.Range("E36").Value = "  -0.25%  "
Set-TextValue "D37" "2.719"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue "D38" "0.01929"
$ws.Range("E38").Value = "  +0.43%  "
Set-TextValue "D39" "2.795"
$ws.Range("E39").Value = "  +0.06%  "
Set-TextValue "D40" "0.4468"
$ws.Range("E40").Value = "  +0.45%  "
Set-TextValue "D41" "73.22"
$ws.Range("E41").Value = "  +0.69%  "
Set-TextValue "D42" "5.971"
$ws.Range("E42").Value = "  -3.18%  "
Set-TextValue "D43" "0.8542"
$ws.Range("E43").Value = "  +2.27%  "
Set-TextValue "D44" "1.003"
$ws.Range("E44").Value = "  +0.19%  "
Set-TextValue "D45" "1.908"
$ws.Range("E45").Value = "  +0.46%  "
Set-TextValue "D46" "102.40"
$ws.Range("E46").Value = "  +1.74%  "
Set-TextValue "D47" "7.594"
$ws.Range("E47").Value = "  +0.17%  "
Set-TextValue "D48" "9.841"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "1.006.28"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "2.063.58"
$ws.Range("E50").Value = "  +0.43%  "
Set-TextValue "D51" "1.558"
$ws.Range("E51").Value = "  +4.14%  "
